$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the version line on the About sheet (A2)
$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"

# Update the recommended citation line on the About sheet (A6)
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Kostromovskaya Coal Mine, Russia, M1343, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Update the version column (S) for rows 2 through 14 on the Boundaries sheet
for ($row = 2; $row -le 14; $row++) {
    $wsBoundaries.Range("S" + $row).Value = "mines - January 30 (built on " + $newStamp + ")"
}
